$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 950
$ws.Range("B2").Value = 909
$ws.Range("C2").Value = 909
$ws.Range("D2").Value = 909
$ws.Range("E2").Value = 959
$ws.Range("F2").Value = 980
$ws.Range("G2").Value = 948
$ws.Range("H2").Value = 979
